$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Schedule")
$ws2 = $wb.Worksheets.Item("Detailed")

# --- Schedule sheet updates (rows 2-5) ---
$ws1.Range("A2").Value = 46037.02083333334
$ws1.Range("C2").Value = 4.5
$ws1.Range("D2").Value = 17.01
$ws1.Range("E2").Value = 661.83422175
$ws1.Range("F2").Value = 38.90853743386243
$ws1.Range("A3").Value = 46037.29166666666
$ws1.Range("C3").Value = 9
$ws1.Range("D3").Value = 34.02
$ws1.Range("E3").Value = 597.29136675
$ws1.Range("F3").Value = 17.55706545414462
$ws1.Range("A4").Value = 46037.89583333334
$ws1.Range("B4").Value = 46038.16666666666
$ws1.Range("C4").Value = 6.5
$ws1.Range("D4").Value = 24.57
$ws1.Range("E4").Value = 718.032588
$ws1.Range("F4").Value = 29.22395555555556
$ws1.Range("A5").Value = 46038.33333333334
$ws1.Range("C5").Value = 8
$ws1.Range("D5").Value = 30.24
$ws1.Range("E5").Value = 540.6628304999999
$ws1.Range("F5").Value = 17.87906185515873

# --- Detailed sheet updates ---
$ws2.Range("E2").Value = "OFF"
$ws2.Range("E15").Value = "OFF"
$ws2.Range("B41").Value = 124.79767
$ws2.Range("B42").Value = 142.36
$ws2.Range("C43").Value = "historical"
$ws2.Range("C44").Value = "historical"
$ws2.Range("B45").Value = 71.4
$ws2.Range("E45").Value = "ON"
$ws2.Range("B46").Value = 84.50611
$ws2.Range("E46").Value = "ON"
$ws2.Range("B47").Value = 74.11643
$ws2.Range("B49").Value = 57.79891
$ws2.Range("B50").Value = 56.98
$ws2.Range("B51").Value = 58.95402
$ws2.Range("B52").Value = 58.95394
$ws2.Range("B53").Value = 58.16581
$ws2.Range("B54").Value = 36.07
$ws2.Range("B55").Value = 50.37846
$ws2.Range("B56").Value = 36.07
$ws2.Range("E56").Value = "ON"
$ws2.Range("B57").Value = 36.07
$ws2.Range("E57").Value = "ON"
$ws2.Range("B59").Value = 64.93029
$ws2.Range("B61").Value = 76.26062
$ws2.Range("B62").Value = 71.95462
$ws2.Range("B64").Value = 40.54
$ws2.Range("E64").Value = "OFF"
$ws2.Range("B65").Value = 41.26969
$ws2.Range("E65").Value = "OFF"
$ws2.Range("B66").Value = 48.26714
$ws2.Range("B67").Value = 42.04025
$ws2.Range("B70").Value = 36.06
$ws2.Range("B71").Value = 36.06
$ws2.Range("B72").Value = 23.10045
$ws2.Range("B73").Value = 36.06
$ws2.Range("B75").Value = 36.0601
$ws2.Range("B77").Value = 36.0601
$ws2.Range("B79").Value = 32.5543
$ws2.Range("B80").Value = 27.01543
$ws2.Range("B81").Value = 20.94801
$ws2.Range("B82").Value = 38.44817
$ws2.Range("B83").Value = 8.27147
$ws2.Range("B84").Value = 11.08967
$ws2.Range("B85").Value = -8.0681
$ws2.Range("B86").Value = -6.78305
$ws2.Range("B87").Value = -3.99001
$ws2.Range("B88").Value = -3.07171
$ws2.Range("B89").Value = 22.01959
$ws2.Range("B90").Value = 29.85322
$ws2.Range("B92").Value = 30.1875
$ws2.Range("B93").Value = 78
$ws2.Range("B94").Value = 64.8901
$ws2.Range("B95").Value = 57.04922
